# "Ricontrollo semifinale del materiale" - small end-to-end proofreading pass.
$p = $ppt.ActivePresentation

# --- Slide 11 ("Media Framework: handle Codec useful for ...") -----------
# Add a comma after "Codec": "handle Codec useful" -> "handle Codec, useful".
# Re-set only the exact run text so the surrounding run is left intact
# (same <a:rPr>, single <a:r>/<a:t> pair, no splitting).
$slide11 = $p.Slides.Item(11)
$contentShape = $slide11.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$fullText = $tr.Text
$needle = ": handle Codec useful for acquisition and playing multimedia contents."
$pos = $fullText.IndexOf($needle)
if ($pos -ge 0) {
    $run = $tr.Characters($pos + 1, $needle.Length)
    $run.Text = ": handle Codec, useful for acquisition and playing multimedia contents."
}

# --- Slide 15 ("Preserve resources: ...") ---------------------------------
# Collapse the three runs " " + "android" + " " into a single, capitalised
# " Android " run (matching the casing/style used elsewhere on the slide).
$slide15 = $p.Slides.Item(15)
$bodyShape = $slide15.Shapes.Item(2)
$tr2 = $bodyShape.TextFrame.TextRange
$fullText2 = $tr2.Text
$search_from = 0
$target = "android"
$foundPos = -1
while ($true) {
    $i = $fullText2.IndexOf($target, $search_from)
    if ($i -lt 0) { break }
    # make sure it is the standalone word " android " (surrounded by spaces),
    # i.e. the occurrence inside "... is necessary android will recreate ..."
    if ($i -gt 0 -and $fullText2.Substring($i - 1, 1) -eq " " -and ($i + $target.Length) -lt $fullText2.Length -and $fullText2.Substring($i + $target.Length, 1) -eq " ") {
        $before = $fullText2.Substring([Math]::Max(0, $i - 12), 12)
        if ($before.Contains("necessary")) {
            $foundPos = $i
            break
        }
    }
    $search_from = $i + 1
}
if ($foundPos -ge 0) {
    $start = $foundPos  # 1-based index of the space right before "android"
    $run2 = $tr2.Characters($start, 9)  # " android " -> 9 characters
    $run2.Text = " Android "
}
